$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing table (rows 2-15 hold index 1-14)
$newRows = @(
    @{ Idx = 15; Text = "contextText Menu strip Control" },
    @{ Idx = 16; Text = "Mdi Form With menu strip Control" },
    @{ Idx = 17; Text = "Tooltip Control" },
    @{ Idx = 18; Text = "Common Dialog Control" },
    @{ Idx = 19; Text = "Rich Text Box Control" },
    @{ Idx = 20; Text = "File Open & Save" }
)

$startRow = 16
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $prevRow = $row - 1

    # Copy formatting from the row directly above so the new row matches
    # the existing style (centered index column, left-aligned text column).
    $ws.Range("A$prevRow`:B$prevRow").Copy()
    $ws.Range("A$row`:B$row").PasteSpecial(-4122)

    $ws.Range("A$row").Value = $newRows[$i].Idx
    $ws.Range("B$row").Value = $newRows[$i].Text
}

$excel.CutCopyMode = 0

$lastRow = $startRow + $newRows.Count - 1

# Update the active selection / view to mirror where Excel would leave the
# cursor after entering the data (last cell of the new range).
$ws.Range("B$lastRow").Select()
$ws.Application.ActiveWindow.ScrollRow = 12
